$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values in columns A and B for rows 1-5
$ws.Range("A1").Value = 0.0009336
$ws.Range("B1").Value = 1.51

$ws.Range("A2").Value = 0.0006491
$ws.Range("B2").Value = 1.335

$ws.Range("A3").Value = 0.0003856
$ws.Range("B3").Value = 1.1

$ws.Range("A4").Value = 0.0001631
$ws.Range("B4").Value = 0.775

$ws.Range("A5").Value = 0.00001835
$ws.Range("B5").Value = 0.27

# Add new row 6 with average formula in C6
$ws.Range("C6").Formula = "=AVERAGE(C1:C5)"

# Update selection to reflect new active cell
$ws.Range("C6").Select()
